$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: shift the quarter labels in column B
#    up by one slot (row2 now shows the newest quarter, 2022-Q4) and append a
#    new trailing row for the quarter that fell off the bottom (2021-Q1).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("B7").Value = "2021-Q3"

$summary.Range("A7").Copy($summary.Range("A8"))
$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 1
$summary.Range("D8").Value = 0.02

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. right before
#    the existing "2022-Q3" sheet), carrying over the same layout/formatting
#    as the other quarterly sheets by copying "2022-Q3" and overwriting the
#    figures with the new quarter's numbers.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.64"
$q4.Range("E2").Value = "93.56"
$q4.Range("F2").Value = "3.19"
$q4.Range("G2").Value = "0.0204"
$q4.Range("H2").Value = 9
